$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style from the last existing date cell (A357) so the new date
# cells in column A get the same formatting (s="2") without introducing
# any new style/font entries in styles.xml.
$ws.Cells.Item(357, 1).Copy()

$newData = @(
    @(358, 44432, 2, 5, 80.11536612722321),
    @(359, 44433, 0, 5, 80.11536612722321),
    @(360, 44434, 1, 6, 96.13843935266785),
    @(361, 44435, 0, 5, 80.11536612722321),
    @(362, 44436, 1, 6, 96.13843935266785),
    @(363, 44437, 0, 4, 64.09229290177856),
    @(364, 44438, 0, 4, 64.09229290177856),
    @(365, 44439, 0, 2, 32.04614645088928),
    @(366, 44440, 0, 2, 32.04614645088928)
)

foreach ($entry in $newData) {
    $r = $entry[0]
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)
    $ws.Cells.Item($r, 1).Value = $entry[1]
    $ws.Cells.Item($r, 2).Value = $entry[2]
    $ws.Cells.Item($r, 3).Value = $entry[3]
    $ws.Cells.Item($r, 4).Value = $entry[4]
}
